$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 330, shifting existing rows 330-424 down to 331-425.
$ws.Rows(330).Insert()

# Populate the newly inserted row 330 with the new weekly record.
$ws.Cells.Item(330, 1).Value = 7
$ws.Cells.Item(330, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(330, 3).Value = "Ñuble"
$ws.Cells.Item(330, 4).Value = 44932
$ws.Cells.Item(330, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(330, 5).Value = 16
$ws.Cells.Item(330, 6).Value = "Fruta"
$ws.Cells.Item(330, 7).Value = 100101
$ws.Cells.Item(330, 8).Value = "Berries"
$ws.Cells.Item(330, 9).Value = 100112025
$ws.Cells.Item(330, 10).Value = "Frutilla"
$ws.Cells.Item(330, 11).Value = "Sin especificar"
$ws.Cells.Item(330, 12).Value = "Especial"
$ws.Cells.Item(330, 13).Value = 80
$ws.Cells.Item(330, 14).Value = 8000
$ws.Cells.Item(330, 15).Value = 8000
$ws.Cells.Item(330, 16).Value = 8000
$ws.Cells.Item(330, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(330, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(330, 19).Value = 1143
$ws.Cells.Item(330, 20).Value = 7
